# Insert a new weekly price record for "Arveja Verde" (Macroferia Regional
# de Talca) as row 134, pushing the existing rows 134:164 down to 135:165.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 134 (inherits formatting, e.g.
# the date number format on column D, from the row above - same as Excel).
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A134").Value = 5
$ws.Range("B134").Value = "Macroferia Regional de Talca"
$ws.Range("C134").Value = "Maule"
$ws.Range("D134").Value = 45258
$ws.Range("E134").Value = 7
$ws.Range("F134").Value = 100112022
$ws.Range("G134").Value = "Arveja Verde"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 600
$ws.Range("K134").Value = 18000
$ws.Range("L134").Value = 20000
$ws.Range("M134").Value = 19000
$ws.Range("N134").Value = "$/saco 25 kilos"
$ws.Range("O134").Value = "Región del Maule"
$ws.Range("P134").Value = 760
$ws.Range("Q134").Value = 25
$ws.Range("R134").Value = "Hortaliza"
